# Adds V10 (WAV) and V11 (Vehicle Inspection Date) to the "Vehicle" sheet,
# and updates the active-tab/selection bookkeeping to match.

$wb = $excel.ActiveWorkbook

# --- Vehicle sheet: append two new field rows -----------------------------
$ws = $wb.Worksheets.Item("Vehicle")

# Row 12: V10 / WAV
$ws.Cells.Item(12, 1).Value = "V10"
$ws.Cells.Item(12, 2).Value = "WAV"
$ws.Cells.Item(12, 3).Value = "String"
$ws.Cells.Item(12, 4).Value = "Yes"
$ws.Cells.Item(12, 5).Value = "Was this vehicle approved by the City of Chicago as a Wheelchair Accessible Vehicle, with that approval being valid on the last day of the reporting period? Enter Y or N only."

$ws.Range("A12:D12").VerticalAlignment = -4160
$ws.Cells.Item(12, 5).VerticalAlignment = -4160
$ws.Cells.Item(12, 5).WrapText = $true
$ws.Rows.Item(12).RowHeight = 47.25

# Row 13: V11 / Vehicle Inspection Date
$ws.Cells.Item(13, 1).Value = "V11"
$ws.Cells.Item(13, 2).Value = "Vehicle Inspection Date"
$ws.Cells.Item(13, 3).Value = "ISO 8601"
$ws.Cells.Item(13, 4).Value = "Yes"
$ws.Cells.Item(13, 5).Value = "The date (no time) when the vehicle was last inspected by a City of Chicago authorized facility and found to be in compliance with all requirements necessary for operating the vehicle for TNP services."

$ws.Cells.Item(13, 1).VerticalAlignment = -4160
$ws.Cells.Item(13, 2).VerticalAlignment = -4160
$ws.Cells.Item(13, 4).VerticalAlignment = -4160
$ws.Cells.Item(13, 5).VerticalAlignment = -4160
$ws.Cells.Item(13, 5).WrapText = $true
$ws.Rows.Item(13).RowHeight = 47.25

# Move the selection past the newly-added rows.
$ws.Range("E14").Select()

# --- Active tab bookkeeping -------------------------------------------------
# Trip sheet no longer holds the selection/active tab; File Names does instead.
$wsTrip = $wb.Worksheets.Item("Trip")
$wsTrip.Range("C6").Select()

$wsFileNames = $wb.Worksheets.Item("File Names")
$wsFileNames.Activate()
$wsFileNames.Range("A1").Select()
